$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: period changes from 2110 -> 2005, amount 29260 -> 35112 (worker stays SENEN)
$ws.Range("E16").Value = "2005"
$ws.Range("F16").Value = 35112

# Row 17: new data - becomes JOSE GREGORIO ALCAZAR ARRIETA, doc 9101444, period 2005
$ws.Range("C17").Value = "9101444"
$ws.Range("D17").Value = "JOSE GREGORIO ALCAZAR ARRIETA"
$ws.Range("E17").Value = "2005"

# Rows 18-32: SENEN's periods shift up by one (chronological reorder), amount stays 35112
$ws.Range("E18").Value = "2007"
$ws.Range("E19").Value = "2008"
$ws.Range("E20").Value = "2009"
$ws.Range("E21").Value = "2010"
$ws.Range("E22").Value = "2011"
$ws.Range("E23").Value = "2012"
$ws.Range("E24").Value = "2101"
$ws.Range("E25").Value = "2102"
$ws.Range("E26").Value = "2103"
$ws.Range("E27").Value = "2104"
$ws.Range("E28").Value = "2105"
$ws.Range("E29").Value = "2106"
$ws.Range("E30").Value = "2107"
$ws.Range("E31").Value = "2108"
$ws.Range("E32").Value = "2109"

# Row 33: becomes last SENEN row with period 2110, amount reverts to 29260
$ws.Range("C33").Value = "8853287"
$ws.Range("D33").Value = "SENEN MARTINEZ TORREGLOSA"
$ws.Range("E33").Value = "2110"
$ws.Range("F33").Value = 29260
